$wb = $excel.ActiveWorkbook

# Remember the originally active sheet so we can restore it after adding
# the new sheet (Excel normally activates a newly inserted sheet).
$originalActive = $wb.ActiveSheet

# --- Add the new worksheet "FTNC_Demand515" as the last sheet ----------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "FTNC_Demand515"

# --- Header row (row 1) --------------------------------------------------
$ws.Range("B1").Value = "In-vehicle"
$ws.Range("C1").Value = "At-stop"
$ws.Range("D1").Value = "Extra"
$ws.Range("E1").Value = "Tardiness"
$ws.Range("F1").Value = "Total"

# --- Data row (row 2) -----------------------------------------------------
$ws.Range("A2").Value = "FTNC"
$ws.Range("B2").Value = 2151.369560509833
$ws.Range("C2").Value = 12665.43114145906
$ws.Range("D2").Value = 461.8074481186643
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 15278.60815008753

# --- Formatting: bold text, centered/top aligned, thin box border --------
$headerRange = $ws.Range("B1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

$labelCell = $ws.Range("A2")
$labelCell.Font.Bold = $true
$labelCell.HorizontalAlignment = -4108  # xlCenter
$labelCell.VerticalAlignment = -4160    # xlTop
$labelCell.Borders.LineStyle = 1
$labelCell.Borders.Weight = 2

$ws.Range("A1").Select()

# Restore the workbook's originally active sheet / tab.
$originalActive.Activate()
